# Add a header row (id, code, name, description) and a sample data row
# (10, 10, 10, 10) to every worksheet in the workbook, and widen the four
# columns so the new header text has room to breathe.

$wb = $excel.ActiveWorkbook

$headers = @("id", "code", "name", "description")
$widths  = @(10, 20, 30, 40)

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item(1, $col).Value = $headers[$i]
        $ws.Cells.Item(2, $col).Value = 10
        $ws.Columns.Item($col).ColumnWidth = $widths[$i]
    }
}
